$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.139.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.900.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  -0.53%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.54%  "

$ws.Range("E6").Value = "  +1.25%  "

$ws.Range("E7").Value = "  -0.50%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.69"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.06%  "

$ws.Range("E9").Value = "  +1.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.54%  "

$ws.Range("E11").Value = "  +5.01%  "

$ws.Range("E12").Value = "  -1.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.176.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("E16").Value = "  +3.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.918.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.138.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0841"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "252.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.69%  "

$ws.Range("E22").Value = "  +1.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.27%  "

$ws.Range("E24").Value = "  -0.45%  "

$ws.Range("E25").Value = "  +4.85%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.27%  "

$ws.Range("E30").Value = "  -0.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.128.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("E32").Value = "  +10.44%  "

$ws.Range("E33").Value = "  +3.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0595"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.35%  "

$ws.Range("E35").Value = "  +10.77%  "

$ws.Range("E36").Value = "  +3.53%  "

$ws.Range("E37").Value = "  -0.53%  "

$ws.Range("E39").Value = "  +0.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.76%  "

$ws.Range("E42").Value = "  +3.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0668"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.59%  "

$ws.Range("E44").Value = "  +0.81%  "

$ws.Range("E45").Value = "  +1.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.308.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.35%  "

$ws.Range("E47").Value = "  -0.06%  "

$ws.Range("E48").Value = "  -1.66%  "

$ws.Range("E49").Value = "  +1.80%  "

$ws.Range("E50").Value = "  -2.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0766"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.19%  "
